$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.012.79'
$ws.Range('E2').Value = '  +2.04%  '
$ws.Range('D3').Value = '1.670.81'
$ws.Range('E3').Value = '  +2.84%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.04'
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('E6').Value = '  +2.08%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +2.14%  '
$ws.Range('E9').Value = '  +1.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.10'
$ws.Range('E10').Value = '  +4.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0894'
$ws.Range('E11').Value = '  +4.55%  '
$ws.Range('D12').Value = '1.905.66'
$ws.Range('E12').Value = '  +2.74%  '
$ws.Range('D13').Value = '1.663.26'
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('E14').Value = '  +1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '65.75'
$ws.Range('E15').Value = '  +2.78%  '
$ws.Range('E16').Value = '  +1.71%  '
$ws.Range('D17').Value = '27.034.35'
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '235.30'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('E19').Value = '  +1.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.73'
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  +3.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.22'
$ws.Range('E23').Value = '  +1.21%  '
$ws.Range('E24').Value = '  +1.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.50'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.17'
$ws.Range('E26').Value = '  +1.35%  '
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.87'
$ws.Range('E28').Value = '  +1.56%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('E31').Value = '  +1.57%  '
$ws.Range('E32').Value = '  +2.04%  '
$ws.Range('D33').Value = '1.450.83'
$ws.Range('E33').Value = '  -4.32%  '
$ws.Range('E34').Value = '  +5.29%  '
$ws.Range('E35').Value = '  +5.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.41'
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.570'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.893'
$ws.Range('E38').Value = '  +7.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0170'
$ws.Range('E39').Value = '  +1.95%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('E42').Value = '  +11.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.29'
$ws.Range('E43').Value = '  +3.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.76'
$ws.Range('E44').Value = '  +4.63%  '
$ws.Range('D45').Value = '1.814.56'
$ws.Range('E45').Value = '  +2.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.779'
$ws.Range('E46').Value = '  +2.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.34'
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('E48').Value = '  +1.52%  '
$ws.Range('E49').Value = '  +4.14%  '
$ws.Range('E50').Value = '  +1.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.63'
$ws.Range('E51').Value = '  +0.77%  '
